$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 44, pushing existing rows 44:49 down to 45:50
$ws.Rows.Item(44).Insert()

# Fill the new row 44 with the new weekly price entry
$ws.Cells.Item(44, 1).Value = 10
$ws.Cells.Item(44, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(44, 3).Value = "La Araucanía"
$ws.Cells.Item(44, 4).Value = 44449
$ws.Cells.Item(44, 4).NumberFormat = $ws.Cells.Item(45, 4).NumberFormat
$ws.Cells.Item(44, 5).Value = 9
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100108
$ws.Cells.Item(44, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(44, 9).Value = 100108004
$ws.Cells.Item(44, 10).Value = "Papaya"
$ws.Cells.Item(44, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 40
$ws.Cells.Item(44, 14).Value = 20000
$ws.Cells.Item(44, 15).Value = 20000
$ws.Cells.Item(44, 16).Value = 20000
$ws.Cells.Item(44, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(44, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(44, 19).Value = 2000
$ws.Cells.Item(44, 20).Value = 10

